$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Fix Runmode for rows 2 and 3 (previously "N", now "Y")
$ws.Range("D2").Value = "Y"
$ws.Range("D3").Value = "Y"

# Add new test case row 14, matching the look & feel of row 2
$ws.Range("A2:E2").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A14").Value = "IPA112"
$ws.Range("B14").Value = "OBT"
$ws.Range("C14").Value = "Save the technology search data and rerun the saved data"
$ws.Range("D14").Value = "Y"
$ws.Range("E14").Value = "SKIP"

# Update view state to match target (scrolled down one row, new selection)
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("D18").Select()
